$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Teste / luisclr.contato@gmail.com
$ws.Range("A3").Value = "Teste"
$ws.Range("B3").Value = "luisclr.contato@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:luisclr.contato@gmail.com") | Out-Null
$ws.Range("B3").Style = "Hiperlink"

# Row 4: email typed first, then name -> contato.jordaquino@gmail.com / Jordan Aquino
$ws.Range("B4").Value = "contato.jordaquino@gmail.com"
$ws.Range("A4").Value = "Jordan Aquino"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:contato.jordaquino@gmail.com") | Out-Null
$ws.Range("B4").Style = "Hiperlink"

# Row 7: a lone formatted (underlined) empty cell
$ws.Range("B7").Font.Underline = 1

$ws.Range("A4").Select()
